$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.250.74"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "1.631.73"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("E6").Value = "  +1.22%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").Value = "1.628.44"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.74%  "

$ws.Range("D16").Value = "27.211.73"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.23%  "

$ws.Range("E23").Value = "  -1.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.118"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("D33").Value = "1.318.71"
$ws.Range("E33").Value = "  +4.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.543"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.849"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.73%  "

$ws.Range("D43").Value = "1.769.63"
$ws.Range("E43").Value = "  -1.27%  "

$ws.Range("E44").Value = "  -4.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "

$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.814"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +21.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.80%  "
